$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column A entries to be stored as literal text (matching the existing
# "Date" column cells, which are shared strings, not real date values) instead
# of letting Excel auto-convert the "2024-10-0x" strings into date serials.
$ws.Range("A24:A25").NumberFormat = "@"

$ws.Range("A24").Value = "2024-10-04"
$ws.Range("B24").Value = 60746.35
$ws.Range("A25").Value = "2024-10-03"
$ws.Range("B25").Value = 60641.8

# Restore the default cell style so the new cells don't end up with a stray
# style index (keeps them identical in appearance to the other data rows).
$ws.Range("A24:A25").Style = "Normal"
